# This script reproduces the cryptocurrency data refresh described in the commit,
# updating Price (column D) and Volume(1h) (column E) figures for each coin row,
# and replacing the "Cronos" row with "Algorand" data in row 51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.260.20"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.57"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("E5").Value = "  +0.78%  "

# Row 6
$ws.Range("E6").Value = "  +0.05%  "

# Row 7
$ws.Range("E7").Value = "  -0.40%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2837"
$ws.Range("E8").Value = "  +0.77%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06526"
$ws.Range("E9").Value = "  -0.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.48"
$ws.Range("E10").Value = "  +6.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07904"
$ws.Range("E11").Value = "  +1.10%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.19"
$ws.Range("E12").Value = "  +0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.861.71"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.156"
$ws.Range("E14").Value = "  +1.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6791"
$ws.Range("E15").Value = "  +2.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.98"
$ws.Range("E16").Value = "  -1.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.249.90"
$ws.Range("E17").Value = "  -0.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.74"
$ws.Range("E18").Value = "  +9.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.393"
$ws.Range("E20").Value = "  -0.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007314"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.107.98"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.147"
$ws.Range("E24").Value = "  +0.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.15"
$ws.Range("E25").Value = "  -0.58%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.183"
$ws.Range("E26").Value = "  -1.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  +0.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("E28").Value = "  +0.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.386"
$ws.Range("E29").Value = "  +3.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09735"
$ws.Range("E30").Value = "  +1.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.376"
$ws.Range("E31").Value = "  -0.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.478"
$ws.Range("E32").Value = "  +0.53%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.039"
$ws.Range("E33").Value = "  -1.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04721"
$ws.Range("E34").Value = "  +1.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +2.77%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7059"
$ws.Range("E36").Value = "  +1.11%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  +0.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.614"
$ws.Range("E39").Value = "  +4.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.327"
$ws.Range("E40").Value = "  -2.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.95"
$ws.Range("E41").Value = "  +3.80%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.955"
$ws.Range("E42").Value = "  +1.10%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8504"
$ws.Range("E43").Value = "  -0.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4169"
$ws.Range("E44").Value = "  +0.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("E46").Value = "  -0.55%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "983.26"
$ws.Range("E47").Value = "  -2.33%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.314"
$ws.Range("E48").Value = "  +2.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.165"
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.06"
$ws.Range("E50").Value = "  +0.83%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1126"
$ws.Range("E51").Value = "  -1.23%  "
